$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure numeric-looking text columns (Price, Volume%, Hora) keep their
# original text representation (leading/trailing zeros, "%", etc.) instead
# of being auto-converted to numbers by Excel when the new value is assigned.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply the refreshed coin data values as per the updated symbol list
$ws.Range("D2").Value = '329.11'
$ws.Range("E2").Value = '0.27%'
$ws.Range("G2").Value = '15'
$ws.Range("D3").Value = '39.73'
$ws.Range("E3").Value = '0.21%'
$ws.Range("G3").Value = '15'
$ws.Range("D4").Value = '5.715'
$ws.Range("E4").Value = '1.43%'
$ws.Range("G4").Value = '15'
$ws.Range("D5").Value = '0.08072'
$ws.Range("E5").Value = '0.13%'
$ws.Range("G5").Value = '15'
$ws.Range("D6").Value = '8.640'
$ws.Range("E6").Value = '-0.59%'
$ws.Range("G6").Value = '15'
$ws.Range("D7").Value = '4.495'
$ws.Range("E7").Value = '-1.46%'
$ws.Range("G7").Value = '15'
$ws.Range("D8").Value = '1.955'
$ws.Range("E8").Value = '0.20%'
$ws.Range("G8").Value = '15'
$ws.Range("G9").Value = '15'
$ws.Range("D10").Value = '0.9271'
$ws.Range("E10").Value = '-1.92%'
$ws.Range("G10").Value = '15'
$ws.Range("D11").Value = '0.1253'
$ws.Range("E11").Value = '-1.92%'
$ws.Range("G11").Value = '15'
$ws.Range("D12").Value = '0.1949'
$ws.Range("E12").Value = '-1.58%'
$ws.Range("G12").Value = '15'
$ws.Range("D13").Value = '8.688'
$ws.Range("E13").Value = '13.48%'
$ws.Range("G13").Value = '15'
$ws.Range("D14").Value = '0.09177'
$ws.Range("E14").Value = '-0.55%'
$ws.Range("G14").Value = '15'
$ws.Range("D15").Value = '0.03641'
$ws.Range("E15").Value = '2.78%'
$ws.Range("G15").Value = '15'
$ws.Range("D16").Value = '0.1053'
$ws.Range("E16").Value = '9.79%'
$ws.Range("G16").Value = '15'
$ws.Range("D17").Value = '0.001300'
$ws.Range("E17").Value = '-1.89%'
$ws.Range("G17").Value = '15'
$ws.Range("D18").Value = '0.006211'
$ws.Range("E18").Value = '-1.66%'
$ws.Range("G18").Value = '15'
$ws.Range("D19").Value = '3.373'
$ws.Range("E19").Value = '0.00%'
$ws.Range("G19").Value = '15'
$ws.Range("E20").Value = '-1.19%'
$ws.Range("G20").Value = '15'
$ws.Range("D21").Value = '0.1367'
$ws.Range("E21").Value = '-2.88%'
$ws.Range("G21").Value = '15'
$ws.Range("D22").Value = '0.2604'
$ws.Range("E22").Value = '3.39%'
$ws.Range("G22").Value = '15'
$ws.Range("D23").Value = '0.04429'
$ws.Range("E23").Value = '-0.49%'
$ws.Range("G23").Value = '15'
$ws.Range("D24").Value = '0.001254'
$ws.Range("E24").Value = '0.02%'
$ws.Range("G24").Value = '15'
$ws.Range("D25").Value = '0.004452'
$ws.Range("E25").Value = '3.60%'
$ws.Range("G25").Value = '15'
$ws.Range("D26").Value = '0.0001241'
$ws.Range("E26").Value = '4.16%'
$ws.Range("G26").Value = '15'
$ws.Range("G27").Value = '15'
$ws.Range("G28").Value = '15'
$ws.Range("G29").Value = '15'
$ws.Range("G30").Value = '15'
$ws.Range("G31").Value = '15'
$ws.Range("G32").Value = '15'
$ws.Range("G33").Value = '15'
$ws.Range("G34").Value = '15'
$ws.Range("G35").Value = '15'
$ws.Range("G36").Value = '15'
$ws.Range("G37").Value = '15'
$ws.Range("G38").Value = '15'
$ws.Range("D39").Value = '0.02704'
$ws.Range("E39").Value = '7.66%'
$ws.Range("G39").Value = '15'
$ws.Range("D40").Value = '0.05481'
$ws.Range("E40").Value = '5.46%'
$ws.Range("G40").Value = '15'
$ws.Range("D41").Value = '0.007533'
$ws.Range("E41").Value = '3.45%'
$ws.Range("G41").Value = '15'
$ws.Range("D42").Value = '0.009834'
$ws.Range("E42").Value = '10.09%'
$ws.Range("G42").Value = '15'
$ws.Range("D43").Value = '0.1418'
$ws.Range("E43").Value = '-0.12%'
$ws.Range("G43").Value = '15'
$ws.Range("D44").Value = '0.002107'
$ws.Range("E44").Value = '-3.69%'
$ws.Range("G44").Value = '15'
$ws.Range("D45").Value = '0.01158'
$ws.Range("E45").Value = '7.16%'
$ws.Range("G45").Value = '15'
$ws.Range("D46").Value = '0.00006773'
$ws.Range("E46").Value = '0.64%'
$ws.Range("G46").Value = '15'
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").Value = '-0.06%'
$ws.Range("G47").Value = '15'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = '0.003062'
$ws.Range("E48").Value = '6.46%'
$ws.Range("G48").Value = '15'
$ws.Range("B49").Value = 'CoinbaseStockToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D49").Value = '0.002274'
$ws.Range("E49").Value = '26.15%'
$ws.Range("G49").Value = '15'
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").Value = '-0.06%'
$ws.Range("G50").Value = '15'
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").Value = '-0.06%'
$ws.Range("G51").Value = '15'
